# Insert one new data row for 2026/01/24 (3rd entry of that date) right
# before the existing row 687 ("2026/12/29"), pushing everything below it
# down by one row. All subsequent rows keep their original values, just
# shifted to row+1, and the sheet grows from D728 to D729.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(687).Insert()

# Leading apostrophe forces the date-looking string to be stored as text
# (matches the existing "日付" column, which is plain text, not a real date).
$ws.Cells.Item(687, 1).Value = "'2026/01/24"
$ws.Cells.Item(687, 2).Value = "土"
$ws.Cells.Item(687, 3).Value = 13
$ws.Cells.Item(687, 4).Value = 13

# Drop the quote-prefix formatting picked up from the apostrophe entry so
# the new cell keeps the sheet's default (unstyled) look, same as its
# neighbours.
$ws.Cells.Item(687, 1).Style = "Normal"
